# Header updates for summer uploads
# Update the header-row labels in Sheet1 to the new wording used for the
# summer reporting upload. Order below matters: it controls the order in
# which the new shared-string entries are appended when the workbook is
# re-saved, which mirrors how the source file was produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Birthdate"
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"
$ws.Range("G1").Value = "Ministry Course Code and Level"

# Make the header row a bit taller to fit the new wrapped text.
$ws.Rows.Item(1).RowHeight = 48

# Leave the cursor/selection where the editor ended up after the change.
$ws.Range("C12").Select()
